# Fruta / hortaliza, semanal
#
# A new weekly price-report record is inserted at row 31 (pushing the
# existing rows 31-57 down to 32-58), for "Arveja Verde" at the
# "Vega Modelo de Temuco" market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31:57 down by one row to make room for the new record.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 10
$ws.Range("B31").Value = "Vega Modelo de Temuco"
$ws.Range("C31").Value = "La Araucanía"
$ws.Range("D31").Value = 44512
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 100112022
$ws.Range("G31").Value = "Arveja Verde"
$ws.Range("H31").Value = "Perfection"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 50
$ws.Range("K31").Value = 16000
$ws.Range("L31").Value = 16000
$ws.Range("M31").Value = 16000
$ws.Range("N31").Value = "$/saco 25 kilos"
$ws.Range("O31").Value = "Región del Maule"
$ws.Range("P31").Value = 640
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
